# This script reproduces the "Fruta / hortaliza, semanal" weekly update for the
# Hortaliza / Betarraga (Terminal La Palmera de La Serena) price history sheet.
#
# The sheet lists weekly price observations (two rows per week: "Primera" and
# "Segunda" quality grades), ordered from most-recent week (row 296) to the
# oldest week (row 423). The update inserts one brand-new week of data at the
# top of the block (rows 296-297) and, as a consequence, pushes every older week
# down by one position; the two weeks that used to be last (422-423) therefore
# reappear once more, now as trailing rows 424-425.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: duplicate the current last two rows (422, 423) into new rows 424-425
# first, before any values are overwritten, so the copy (values + number
# formatting, e.g. the date format on column D) matches the originals exactly.
$src = $ws.Range("A422:R423")
$dst = $ws.Range("A424:R425")
$src.Copy($dst)

# Step 2: shift the price history down by one week for rows 296-423, i.e. each
# row takes on the Fecha/Volumen/Precio values that used to belong to the row
# two positions above it, except for the first week (rows 296-297) which holds
# the brand-new data point being added.
# Columns: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
#          M=Precio promedio ponderado, P=Precio $/Kg
$data = @(
    ,@(296, 45006, 2000, 500, 600, 550, 183)
    ,@(297, 45006, 1100, 400, 450, 425, 142)
    ,@(298, 44959, 2200, 500, 600, 550, 183)
    ,@(299, 44959, 1360, 400, 450, 425, 142)
    ,@(300, 44196, 3000, 400, 500, 450, 150)
    ,@(301, 44196, 1600, 300, 350, 325, 108)
    ,@(302, 44376, 3400, 500, 550, 525, 175)
    ,@(303, 44376, 1640, 350, 400, 375, 125)
    ,@(304, 44294, 3300, 450, 500, 475, 158)
    ,@(305, 44294, 1540, 350, 400, 375, 125)
    ,@(306, 44551, 2800, 450, 500, 475, 158)
    ,@(307, 44551, 1400, 350, 400, 375, 125)
    ,@(308, 44166, 2700, 400, 500, 450, 150)
    ,@(309, 44166, 1600, 300, 350, 325, 108)
    ,@(310, 44798, 2000, 600, 700, 650, 217)
    ,@(311, 44798, 1520, 500, 550, 525, 175)
    ,@(312, 44691, 2480, 500, 600, 550, 183)
    ,@(313, 44691, 1540, 400, 450, 425, 142)
    ,@(314, 44943, 2000, 500, 600, 550, 183)
    ,@(315, 44943, 1560, 400, 450, 425, 142)
    ,@(316, 44420, 3500, 450, 500, 475, 158)
    ,@(317, 44420, 1600, 350, 400, 375, 125)
    ,@(318, 44637, 2200, 500, 600, 550, 183)
    ,@(319, 44637, 1440, 400, 450, 425, 142)
    ,@(320, 44922, 2200, 550, 600, 575, 192)
    ,@(321, 44922, 1540, 450, 500, 475, 158)
    ,@(322, 44222, 2800, 450, 500, 475, 158)
    ,@(323, 44222, 1600, 350, 400, 375, 125)
    ,@(324, 44343, 3500, 500, 600, 550, 183)
    ,@(325, 44343, 1660, 400, 450, 425, 142)
    ,@(326, 44350, 3500, 500, 600, 550, 183)
    ,@(327, 44350, 1660, 400, 450, 425, 142)
    ,@(328, 44721, 2500, 500, 600, 550, 183)
    ,@(329, 44721, 1540, 400, 450, 425, 142)
    ,@(330, 44719, 2500, 500, 600, 550, 183)
    ,@(331, 44719, 1520, 400, 450, 425, 142)
    ,@(332, 44658, 2200, 450, 500, 475, 158)
    ,@(333, 44658, 1460, 350, 400, 375, 125)
    ,@(334, 44399, 3460, 450, 500, 475, 158)
    ,@(335, 44399, 1600, 350, 400, 375, 125)
    ,@(336, 44273, 3200, 450, 500, 475, 158)
    ,@(337, 44273, 1500, 350, 400, 375, 125)
    ,@(338, 44455, 3400, 450, 500, 475, 158)
    ,@(339, 44455, 1600, 350, 400, 375, 125)
    ,@(340, 44782, 2440, 600, 700, 650, 217)
    ,@(341, 44782, 1520, 500, 550, 525, 175)
    ,@(342, 44749, 2440, 500, 600, 550, 183)
    ,@(343, 44749, 1540, 400, 450, 425, 142)
    ,@(344, 44504, 3000, 450, 500, 475, 158)
    ,@(345, 44504, 1460, 350, 400, 375, 125)
    ,@(346, 44194, 2700, 400, 500, 450, 150)
    ,@(347, 44194, 1540, 300, 350, 325, 108)
    ,@(348, 44435, 6800, 450, 500, 475, 158)
    ,@(349, 44435, 3200, 350, 400, 375, 125)
    ,@(350, 44425, 3400, 450, 500, 475, 158)
    ,@(351, 44425, 1600, 350, 400, 375, 125)
    ,@(352, 44987, 2000, 500, 600, 550, 183)
    ,@(353, 44987, 1460, 400, 450, 425, 142)
    ,@(354, 44371, 3400, 450, 500, 475, 158)
    ,@(355, 44371, 1600, 350, 400, 375, 125)
    ,@(356, 44392, 3400, 450, 500, 475, 158)
    ,@(357, 44392, 1600, 350, 400, 375, 125)
    ,@(358, 44586, 3100, 450, 500, 475, 158)
    ,@(359, 44586, 1520, 350, 400, 375, 125)
    ,@(360, 44901, 2000, 550, 600, 575, 192)
    ,@(361, 44901, 1460, 450, 500, 475, 158)
    ,@(362, 44411, 3400, 450, 500, 475, 158)
    ,@(363, 44411, 1600, 350, 400, 375, 125)
    ,@(364, 44789, 2400, 600, 700, 650, 217)
    ,@(365, 44789, 1540, 500, 550, 525, 175)
    ,@(366, 44859, 2000, 550, 600, 575, 192)
    ,@(367, 44859, 1400, 450, 500, 475, 158)
    ,@(368, 44383, 3400, 450, 500, 475, 158)
    ,@(369, 44383, 1600, 350, 400, 375, 125)
    ,@(370, 44873, 2000, 650, 700, 675, 225)
    ,@(371, 44873, 1500, 550, 600, 575, 192)
    ,@(372, 44238, 3200, 450, 500, 475, 158)
    ,@(373, 44238, 1600, 350, 400, 375, 125)
    ,@(374, 44229, 2800, 450, 500, 475, 158)
    ,@(375, 44229, 1600, 350, 400, 375, 125)
    ,@(376, 44663, 2400, 450, 500, 475, 158)
    ,@(377, 44663, 1560, 350, 400, 375, 125)
    ,@(378, 44572, 3160, 450, 500, 475, 158)
    ,@(379, 44572, 1560, 350, 400, 375, 125)
    ,@(380, 44938, 2000, 500, 600, 550, 183)
    ,@(381, 44938, 1540, 400, 450, 425, 142)
    ,@(382, 44285, 3100, 450, 500, 475, 158)
    ,@(383, 44285, 1680, 350, 400, 375, 125)
    ,@(384, 44292, 3200, 450, 500, 475, 158)
    ,@(385, 44292, 1660, 350, 400, 375, 125)
    ,@(386, 44868, 2200, 550, 600, 575, 192)
    ,@(387, 44868, 1500, 450, 500, 475, 158)
    ,@(388, 44740, 2440, 600, 700, 650, 217)
    ,@(389, 44740, 1400, 500, 550, 525, 175)
    ,@(390, 44910, 2000, 550, 600, 575, 192)
    ,@(391, 44910, 1540, 450, 500, 475, 158)
    ,@(392, 44278, 3000, 450, 500, 475, 158)
    ,@(393, 44278, 1600, 350, 400, 375, 125)
    ,@(394, 44957, 1800, 500, 600, 550, 183)
    ,@(395, 44957, 1400, 400, 450, 425, 142)
    ,@(396, 44651, 2200, 450, 500, 475, 158)
    ,@(397, 44651, 1460, 350, 400, 375, 125)
    ,@(398, 44530, 3000, 450, 500, 475, 158)
    ,@(399, 44530, 1500, 350, 400, 375, 125)
    ,@(400, 44306, 3200, 550, 600, 575, 192)
    ,@(401, 44306, 1700, 450, 500, 475, 158)
    ,@(402, 44413, 3500, 450, 500, 475, 158)
    ,@(403, 44413, 1600, 350, 400, 375, 125)
    ,@(404, 44257, 3200, 450, 500, 475, 158)
    ,@(405, 44257, 1680, 350, 400, 375, 125)
    ,@(406, 44567, 3060, 450, 500, 475, 158)
    ,@(407, 44567, 1520, 350, 400, 375, 125)
    ,@(408, 44364, 3400, 500, 600, 550, 183)
    ,@(409, 44364, 1600, 400, 450, 425, 142)
    ,@(410, 44215, 2800, 450, 500, 475, 158)
    ,@(411, 44215, 1560, 350, 400, 375, 125)
    ,@(412, 44168, 2800, 400, 500, 450, 150)
    ,@(413, 44168, 1600, 300, 350, 325, 108)
    ,@(414, 44677, 2400, 450, 500, 475, 158)
    ,@(415, 44677, 1500, 350, 400, 375, 125)
    ,@(416, 44747, 2440, 500, 600, 550, 183)
    ,@(417, 44747, 1560, 400, 450, 425, 142)
    ,@(418, 44245, 3200, 450, 500, 475, 158)
    ,@(419, 44245, 1600, 350, 400, 375, 125)
    ,@(420, 44236, 2800, 450, 500, 475, 158)
    ,@(421, 44236, 1560, 350, 400, 375, 125)
    ,@(422, 44210, 3000, 450, 500, 475, 158)
    ,@(423, 44210, 1600, 350, 400, 375, 125)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 16).Value = $row[6]
}
